# Auto-generated Excel COM-interop script
# Applies Bahamut_Profits leve-profit recalculations scraped from the scheduled runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: 'Stuck in the Moment' / 'Horn Glue'
$ws.Range("H40").Value = 2447.1667
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2447.1667
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2447.1667
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -2797.1667

# Row 76: 'Warding Off Temptation' / 'Enchanted Hardsilver Ink'
$ws.Range("H76").Value = 42310550
$ws.Range("I76").Value = 44002850
$ws.Range("J76").Value = 3200
$ws.Range("K76").Value = 44002850
$ws.Range("L76").Value = 3200
$ws.Range("M76").Value = -44002535
$ws.Range("N76").Value = -3830

# Row 79: 'The Garden of Arcane Delights (L)' / 'Enchanted Hardsilver Ink'
$ws.Range("H79").Value = 42310550
$ws.Range("I79").Value = 44002850
$ws.Range("J79").Value = 3200
$ws.Range("K79").Value = 44002850
$ws.Range("L79").Value = 3200
$ws.Range("M79").Value = -44001758
$ws.Range("N79").Value = -5384

# Row 118: 'Crafty Concoctions' / "Commanding Craftsman's Syrup"
$ws.Range("H118").Value = 452
$ws.Range("I118").Value = 315
$ws.Range("K118").Value = 945
$ws.Range("M118").Value = 712

# Row 132: 'Fast-forwarding Flora' / 'Growth Formula Lambda'
$ws.Range("H132").Value = 305050.8
$ws.Range("I132").Value = 2052.6667
$ws.Range("J132").Value = 1113045.9
$ws.Range("K132").Value = 6158.000100000001
$ws.Range("L132").Value = 3339137.7
$ws.Range("M132").Value = -3628.000100000001
$ws.Range("N132").Value = -3344197.7

# Row 138: 'All-night Crafting' / "Cunning Craftsman's Tisane"
$ws.Range("H138").Value = 3501.7
$ws.Range("J138").Value = 4051.2651
$ws.Range("L138").Value = 12153.7953
$ws.Range("N138").Value = -22433.7953

$ws = $wb.Worksheets.Item("ARM")
# Row 32: 'Ingot We Trust' / 'Steel Ingot'
$ws.Range("H32").Value = 2279150.2
$ws.Range("I32").Value = 5515.9565
$ws.Range("K32").Value = 5515.9565
$ws.Range("M32").Value = -5228.9565

# Row 128: 'Heading toward Bankruptcy' / 'Manganese Helm of the Falling Dragon'
$ws.Range("H128").Value = 47986.668
$ws.Range("J128").Value = 47986.668
$ws.Range("L128").Value = 47986.668
$ws.Range("N128").Value = -57946.668

# Row 129: 'In-kweh-dible Cooking' / 'Manganese Chocobo Frypan'
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999

$ws = $wb.Worksheets.Item("BSM")
# Row 69: 'Freight and Barrel' / 'Mythrite Stilettos'
$ws.Range("H69").Value = 40000
$ws.Range("J69").Value = 40000
$ws.Range("L69").Value = 40000
$ws.Range("N69").Value = -41622

# Row 72: 'Stepping on My Heart with Stilettos (L)' / 'Mythrite Stilettos'
$ws.Range("H72").Value = 40000
$ws.Range("J72").Value = 40000
$ws.Range("L72").Value = 120000
$ws.Range("N72").Value = -128112

# Row 75: 'I Saw the Pine' / 'Hardsilver Saw'
$ws.Range("H75").Value = 40000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 40000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 40000
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -41872

# Row 76: 'Keep Up with the Mechanics' / 'Titanium-barreled Arquebus'
$ws.Range("H76").Value = 31657
$ws.Range("J76").Value = 31657
$ws.Range("L76").Value = 31657
$ws.Range("N76").Value = -32287

# Row 78: 'I Came, I Sawed, I Conquered (L)' / 'Hardsilver Saw'
$ws.Range("H78").Value = 40000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 40000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 120000
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -129360

# Row 79: 'Unconventional Weaponry (L)' / 'Titanium-barreled Arquebus'
$ws.Range("H79").Value = 31657
$ws.Range("J79").Value = 31657
$ws.Range("L79").Value = 31657
$ws.Range("N79").Value = -33841

# Row 82: 'Spirituality Inspector' / 'Titanium Lump Hammer'
$ws.Range("H82").Value = 17856.125
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

# Row 85: 'The Clamor for Hammers (L)' / 'Titanium Lump Hammer'
$ws.Range("H85").Value = 17856.125
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

# Row 88: 'Swords for Plowshares' / 'Adamantite Zweihander'
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

# Row 91: 'Negative, They Are Meat Popsicles (L)' / 'Adamantite Zweihander'
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

# Row 105: 'Ingot to Wing It' / 'Molybdenum Ingot'
$ws.Range("H105").Value = 4327.1333
$ws.Range("I105").Value = 3923.32
$ws.Range("J105").Value = 6346.2
$ws.Range("K105").Value = 3923.32
$ws.Range("L105").Value = 6346.2
$ws.Range("M105").Value = -2176.32
$ws.Range("N105").Value = -9840.200000000001

# Row 129: 'Pruned to Perfection' / 'Manganese Garden Scythe'
$ws.Range("H129").Value = 49992.668
$ws.Range("J129").Value = 49992.668
$ws.Range("L129").Value = 49992.668
$ws.Range("N129").Value = -59992.668

# Row 134: 'Ruthenium Supremium' / 'Ruthenium Ingot'
$ws.Range("H134").Value = 59141.11
$ws.Range("I134").Value = 3509.6924
$ws.Range("J134").Value = 203782.8
$ws.Range("K134").Value = 10529.0772
$ws.Range("L134").Value = 611348.3999999999
$ws.Range("M134").Value = -7994.0772
$ws.Range("N134").Value = -616418.3999999999

# Row 137: 'Dagger Swagger' / 'Cobalt Tungsten Khukuri'
$ws.Range("H137").Value = 53030.77
$ws.Range("J137").Value = 53030.77
$ws.Range("L137").Value = 53030.77
$ws.Range("N137").Value = -63230.77

$ws = $wb.Worksheets.Item("CRP")
# Row 132: 'Hull Lotta Damage' / 'Ginseng Lumber'
$ws.Range("H132").Value = 2559.353
$ws.Range("I132").Value = 1504.75
$ws.Range("K132").Value = 4514.25
$ws.Range("M132").Value = -1984.25

# Row 140: 'Spear Pressure' / 'Claro Walnut Spear'
$ws.Range("H140").Value = 52830.91
$ws.Range("J140").Value = 52830.91
$ws.Range("L140").Value = 52830.91
$ws.Range("N140").Value = -63190.91

$ws = $wb.Worksheets.Item("CUL")
# Row 5: 'What a Sap' / 'Maple Syrup'
$ws.Range("H5").Value = 1850.091
$ws.Range("I5").Value = 1532.1765
$ws.Range("J5").Value = 2931
$ws.Range("K5").Value = 4596.529500000001
$ws.Range("L5").Value = 8793
$ws.Range("M5").Value = -4484.529500000001
$ws.Range("N5").Value = -9017

# Row 23: 'Sweet Smell of Success' / 'Lavender Oil'
$ws.Range("H23").Value = 39.555557
$ws.Range("I23").Value = 47.6
$ws.Range("J23").Value = 29.5
$ws.Range("K23").Value = 142.8
$ws.Range("L23").Value = 88.5
$ws.Range("M23").Value = 92.19999999999999
$ws.Range("N23").Value = -558.5

# Row 114: 'One Last Meal' / 'Mushroom Saute'
$ws.Range("H114").Value = 1458.037
$ws.Range("I114").Value = 747.7143
$ws.Range("J114").Value = 1706.65
$ws.Range("K114").Value = 2243.1429
$ws.Range("L114").Value = 5119.950000000001
$ws.Range("M114").Value = 1010.8571
$ws.Range("N114").Value = -11627.95

# Row 117: 'A Good Omen' / 'Peppered Popotoes'
$ws.Range("H117").Value = 2376.4119
$ws.Range("I117").Value = 952.6667
$ws.Range("J117").Value = 2681.5
$ws.Range("K117").Value = 2858.0001
$ws.Range("L117").Value = 8044.5
$ws.Range("M117").Value = 583.9998999999998
$ws.Range("N117").Value = -14928.5

# Row 122: 'Salt of the North' / 'Northern Sea Salt'
$ws.Range("H122").Value = 435576.1
$ws.Range("I122").Value = 487.6
$ws.Range("J122").Value = 770259.5600000001
$ws.Range("K122").Value = 4388.400000000001
$ws.Range("L122").Value = 6932336.040000001
$ws.Range("M122").Value = -1938.400000000001
$ws.Range("N122").Value = -6937236.040000001

# Row 129: 'Comfort Food' / 'Yakow Moussaka'
$ws.Range("H129").Value = 2099.9412
$ws.Range("I129").Value = 1141.9
$ws.Range("J129").Value = 2499.125
$ws.Range("K129").Value = 3425.7
$ws.Range("L129").Value = 7497.375
$ws.Range("M129").Value = 1574.3
$ws.Range("N129").Value = -17497.375

# Row 135: 'Not-so-secret Ingredient' / 'Royal Maple Syrup'
$ws.Range("H135").Value = 1850.091
$ws.Range("I135").Value = 1532.1765
$ws.Range("J135").Value = 2931
$ws.Range("K135").Value = 13789.5885
$ws.Range("L135").Value = 26379
$ws.Range("M135").Value = -11254.5885
$ws.Range("N135").Value = -31449

$ws = $wb.Worksheets.Item("GSM")
# Row 70: 'Sky Is the Limit' / 'Mythrite Ingot'
$ws.Range("H70").Value = 4237.9414
$ws.Range("I70").Value = 4049.6667
$ws.Range("J70").Value = 5650
$ws.Range("K70").Value = 4049.6667
$ws.Range("L70").Value = 5650
$ws.Range("M70").Value = -3779.6667
$ws.Range("N70").Value = -6190

# Row 73: 'Hulls of Broken Dreams (L)' / 'Mythrite Ingot'
$ws.Range("H73").Value = 4237.9414
$ws.Range("I73").Value = 4049.6667
$ws.Range("J73").Value = 5650
$ws.Range("K73").Value = 4049.6667
$ws.Range("L73").Value = 5650
$ws.Range("M73").Value = -3113.6667
$ws.Range("N73").Value = -7522

# Row 80: 'Needs More Prayerbell' / 'Hardsilver Ingot'
$ws.Range("H80").Value = 3791.6667
$ws.Range("I80").Value = 3772.7273
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 3772.7273
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -2774.7273
$ws.Range("N80").Value = -5996

# Row 83: 'With a Noise That Reaches Heaven (L)' / 'Hardsilver Ingot'
$ws.Range("H83").Value = 3791.6667
$ws.Range("I83").Value = 3772.7273
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 18863.6365
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -13871.6365
$ws.Range("N83").Value = -29984

# Row 128: 'To Fight at Her Side' / 'Manganese Rapier'
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 32: 'Men Who Scare Up Goats' / 'Goatskin Targe'
$ws.Range("H32").Value = 25000450
$ws.Range("I32").Value = 25000450
$ws.Range("K32").Value = 25000450
$ws.Range("M32").Value = -25000133

# Row 82: "Trainin' the Neck" / 'Dragon Leather'
$ws.Range("H82").Value = 4300
$ws.Range("I82").Value = 3683.3333
$ws.Range("J82").Value = 8000
$ws.Range("K82").Value = 3683.3333
$ws.Range("L82").Value = 8000
$ws.Range("M82").Value = -3322.3333
$ws.Range("N82").Value = -8722

# Row 85: 'Training Is Only Skintight (L)' / 'Dragon Leather'
$ws.Range("H85").Value = 4300
$ws.Range("I85").Value = 3683.3333
$ws.Range("J85").Value = 8000
$ws.Range("K85").Value = 3683.3333
$ws.Range("L85").Value = 8000
$ws.Range("M85").Value = -2435.3333
$ws.Range("N85").Value = -10496

# Row 93: 'Hide to Go Seek' / 'Gagana Leather'
$ws.Range("H93").Value = 1994.125
$ws.Range("I93").Value = 1211.4445
$ws.Range("J93").Value = 3000.4285
$ws.Range("K93").Value = 1211.4445
$ws.Range("L93").Value = 3000.4285
$ws.Range("M93").Value = 36.55549999999994
$ws.Range("N93").Value = -5496.4285

# Row 136: "Respect for Br'aax" / "Br'aax Leather"
$ws.Range("H136").Value = 4041
$ws.Range("I136").Value = 2000.7273
$ws.Range("J136").Value = 7247.143
$ws.Range("K136").Value = 6002.1819
$ws.Range("L136").Value = 21741.429
$ws.Range("M136").Value = -3452.1819
$ws.Range("N136").Value = -26841.429

# Row 139: 'Giving Gatherers Their Gear' / 'Gomphotherium Doublet of Gathering'
$ws.Range("H139").Value = 55715
$ws.Range("J139").Value = 55715
$ws.Range("L139").Value = 55715
$ws.Range("N139").Value = -65995

$ws = $wb.Worksheets.Item("WVR")
# Row 132: 'Comfy Cabins' / 'Snow Cotton Cloth'
$ws.Range("H132").Value = 1764.75
$ws.Range("I132").Value = 1510.5
$ws.Range("J132").Value = 2400.375
$ws.Range("K132").Value = 4531.5
$ws.Range("L132").Value = 7201.125
$ws.Range("M132").Value = -2001.5
$ws.Range("N132").Value = -12261.125

# Row 136: 'Weaving the Envelope' / 'Sarcenet Cloth'
$ws.Range("H136").Value = 2751.52
$ws.Range("I136").Value = 2704.1904
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 8112.5712
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -5562.5712
$ws.Range("N136").Value = -14100
